# Applies the "Add chatting, Fix multicast bug" edit to the packet model workbook.
$wb = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item(1)   # "Packet" sheet
$ws2 = $wb.Worksheets.Item(2)   # "Error Code" sheet

# ---------------------------------------------------------------------------
# Row 11 : the "room update (multicast)" block moved from F/I/J into new spots
# ---------------------------------------------------------------------------
$ws.Range("D11").Value = "[68], [70], [72]"
$ws.Range("F11").Value = "room update (multicast)"
$ws.Range("I11").Value = "s [room info]"
$ws.Range("J11").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# Row 12 : room chat (multicast) block, shifted one column left, new I/J text
# ---------------------------------------------------------------------------
$ws.Range("C12").Value = "71(G)"
$ws.Range("D12").Value = "room chat R"
$ws.Range("E12").Value = "102(f)"
$ws.Range("F12").Value = "room chat (multicast)"
$ws.Range("H12").Value = "[chat]"
$ws.Range("I12").Value = "[sender]\r\n[chat]"
$ws.Range("J12").Value = "\r\n: separator"

# ---------------------------------------------------------------------------
# Row 13 : room exit block
# ---------------------------------------------------------------------------
$ws.Range("C13").Value = "72(H)"
$ws.Range("D13").Value = "room exit R"
$ws.Range("E13").ClearContents() | Out-Null
$ws.Range("F13").Value = "[101]"
$ws.Range("H13").Value = "[room number]"
$ws.Range("I13").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# Row 14 : now "sign up"
# ---------------------------------------------------------------------------
$ws.Range("C14").Value = "73(I)"
$ws.Range("D14").Value = "sign up R"
$ws.Range("E14").Value = "103(g)"
$ws.Range("F14").Value = "sign up"
$ws.Range("H14").Value = "[id] [password]"
$ws.Range("I14").Value = "s / f"

# ---------------------------------------------------------------------------
# Row 15 : now "show friends"
# ---------------------------------------------------------------------------
$ws.Range("C15").Value = "74(J)"
$ws.Range("D15").Value = "show friends R"
$ws.Range("E15").Value = "104(h)"
$ws.Range("F15").Value = "show friends"
$ws.Range("H15").Value = "<empty>"
$ws.Range("I15").Value = "[friends list]"

# ---------------------------------------------------------------------------
# Row 16 : now "add friend"
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = "75(K)"
$ws.Range("D16").Value = "add friend R"
$ws.Range("E16").Value = "105(i)"
$ws.Range("F16").Value = "add friend"
$ws.Range("H16").Value = "[friend id]"
$ws.Range("I16").Value = "s / f / f1"
$ws.Range("J16").Value = "f1: add myself"

# ---------------------------------------------------------------------------
# Row 17 : now "delete friend"
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = "76(L)"
$ws.Range("D17").Value = "delete friend R"
$ws.Range("E17").ClearContents() | Out-Null
$ws.Range("H17").Value = "[friend id]"

# ---------------------------------------------------------------------------
# Rows 18-28 : byte-code column shifts down by 6 (108..118)
# ---------------------------------------------------------------------------
$ws.Range("E18").Value = "106(j)"
$ws.Range("E19").Value = "107"
$ws.Range("E20").Value = "108"
$ws.Range("E21").Value = "109"
$ws.Range("E22").Value = "110"
$ws.Range("E23").Value = "111"
$ws.Range("E24").Value = "112"
$ws.Range("E25").Value = "113"
$ws.Range("E26").Value = "114"
$ws.Range("E27").Value = "115"
$ws.Range("E28").Value = "116"

# ---------------------------------------------------------------------------
# Row 29
# ---------------------------------------------------------------------------
$ws.Range("C29").Value = "88(X)"
$ws.Range("E29").Value = "117"

# Row 30
$ws.Range("E30").Value = "118"

# Row 31
$ws.Range("E31").Value = "119"

# Row 32 : "unknown error" text removed from F32 (moved to row 34)
$ws.Range("E32").Value = "120"
$ws.Range("F32").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# Row 33 : new "notice" entry (121(x))
# ---------------------------------------------------------------------------
$ws.Range("E33").Value = "121(x)"
$ws.Range("F33").Value = "notice"
$ws.Range("F33").NumberFormat = $ws2.Range("B4").NumberFormat
$ws.Range("F33").HorizontalAlignment = $ws2.Range("B4").HorizontalAlignment
$ws.Range("F33").VerticalAlignment = $ws2.Range("B4").VerticalAlignment

# ---------------------------------------------------------------------------
# Row 34 : "unknown error" now here (122(z))
# ---------------------------------------------------------------------------
$ws.Range("E34").Value = "122(z)"
$ws.Range("F34").Value = "unknown error"

# ---------------------------------------------------------------------------
# Column J got wider to fit the new "\r\n: separator" note
# ---------------------------------------------------------------------------
$ws.Columns.Item(10).ColumnWidth = 14.79

# ---------------------------------------------------------------------------
# Selection moved from H32 to J12, and the frozen/scrolled top-left cell
# reset back to A1 (no more topLeftCell="A10")
# ---------------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A1"), $true) | Out-Null
$ws.Range("J12").Select() | Out-Null
